$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 82) down onto the
# new row 83 first, so the appended cells pick up the same cell styles
# (date number format in column A, centered/plain style in B:F) instead of
# Excel synthesizing brand-new style records for a freshly typed value.
$ws.Range("A82:F82").Copy()
$ws.Range("A83:F83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Enter the new day's figures (03/06/2020).
$ws.Range("A83").Value = 43985
$ws.Range("B83").Value = 549
$ws.Range("C83").Value = 163
$ws.Range("D83").Value = 413
$ws.Range("E83").Value = 92
$ws.Range("F83").Value = 44

# The table "Condicion_Pacientes" grows to include the freshly entered row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F83"))

# Match the author's resulting view state: selection parked on the
# newly entered last cell.
$ws.Range("F83").Select()
